# Framework Update with login module
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (TestCaseID 1) - ActiveFlag column (D2) becomes the text "Yes"
# instead of the numeric value 1.
$ws.Range("D2").Value = "Yes"

# The placeholder second test case row (A3:D3 - "Write few test cases
# below") is removed entirely.
$ws.Range("A3:D3").ClearContents()

# Leave the active selection on the now-empty row below the data,
# spanning the full width of the table (A3:D3).
$ws.Range("A3:D3").Select()
